$wb = $excel.ActiveWorkbook

$desired = $wb.Worksheets.Item("Desired State")
$desired.Delete()

$current = $wb.Worksheets.Item("Current State")
$current.Name = "Power Armor Sets"

Write-Host "Sheets after edit:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
